# chore: update Sheets via scheduled runner
# Refreshes cached market-price figures (columns H-N) across the ALC, ARM,
# BSM, CRP, CUL, LTW and WVR sheets. A handful of rows also gain/lose a
# trailing NQ/HQ profit cell as the source feed's column coverage shifts.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2232.6052
$ws.Range("I19").Value = 3970.625
$ws.Range("J19").Value = 968.5909
$ws.Range("K19").Value = 3970.625
$ws.Range("L19").Value = 968.5909
$ws.Range("M19").Value = -3795.625
$ws.Range("N19").Value = -1318.5909
$ws.Range("H33").Value = 446
$ws.Range("I33").Value = 463
$ws.Range("J33").Value = 55
$ws.Range("K33").Value = 463
$ws.Range("L33").Value = 55
$ws.Range("M33").Value = -234
$ws.Range("N33").Value = -513
$ws.Range("H76").Value = 13164345
$ws.Range("I76").Value = 29420862
$ws.Range("J76").Value = 4307.7144
$ws.Range("K76").Value = 29420862
$ws.Range("L76").Value = 4307.7144
$ws.Range("M76").Value = -29420547
$ws.Range("N76").Value = -4937.7144
$ws.Range("H79").Value = 13164345
$ws.Range("I79").Value = 29420862
$ws.Range("J79").Value = 4307.7144
$ws.Range("K79").Value = 29420862
$ws.Range("L79").Value = 4307.7144
$ws.Range("M79").Value = -29419770
$ws.Range("N79").Value = -6491.7144
$ws.Range("H113").Value = 2857.465
$ws.Range("I113").Value = 1942.2222
$ws.Range("J113").Value = 3516.44
$ws.Range("K113").Value = 1942.2222
$ws.Range("L113").Value = 3516.44
$ws.Range("M113").Value = 1311.7778
$ws.Range("N113").Value = -10024.44
$ws.Range("H116").Value = 5163.1665
$ws.Range("I116").Value = 5274.75
$ws.Range("K116").Value = 5274.75
$ws.Range("M116").Value = -1832.75
$ws.Range("H132").Value = 6084
$ws.Range("I132").Value = 1422
$ws.Range("J132").Value = 14242.5
$ws.Range("K132").Value = 4266
$ws.Range("L132").Value = 42727.5
$ws.Range("M132").Value = -1736
$ws.Range("N132").Value = -47787.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 850.5625
$ws.Range("I2").Value = 850.5625
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 850.5625
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -737.5625
$ws.Range("H45").Value = 979.875
$ws.Range("I45").Value = 797.7778
$ws.Range("J45").Value = 1214
$ws.Range("K45").Value = 797.7778
$ws.Range("L45").Value = 1214
$ws.Range("M45").Value = -420.7778
$ws.Range("N45").Value = -1968
$ws.Range("H61").Value = 292615.2
$ws.Range("I61").Value = 205384.84
$ws.Range("J61").Value = 506329.6
$ws.Range("K61").Value = 205384.84
$ws.Range("L61").Value = 506329.6
$ws.Range("M61").Value = -205172.84
$ws.Range("N61").Value = -506753.6
$ws.Range("H63").Value = 11560
$ws.Range("I63").Value = 17933.334
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 17933.334
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -17247.334
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 11560
$ws.Range("I66").Value = 17933.334
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 89666.67
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -86234.67
$ws.Range("N66").Value = -16864
$ws.Range("H116").Value = 850.5625
$ws.Range("I116").Value = 850.5625
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 850.5625
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1443.4375
$ws.Range("H132").Value = 30185.973
$ws.Range("I132").Value = 38874.32
$ws.Range("J132").Value = 3155.5557
$ws.Range("K132").Value = 116622.96
$ws.Range("L132").Value = 9466.667099999999
$ws.Range("M132").Value = -114092.96
$ws.Range("N132").Value = -14526.6671
$ws.Range("H136").Value = 292615.2
$ws.Range("I136").Value = 205384.84
$ws.Range("J136").Value = 506329.6
$ws.Range("K136").Value = 616154.52
$ws.Range("L136").Value = 1518988.8
$ws.Range("M136").Value = -613604.52
$ws.Range("N136").Value = -1524088.8
$ws.Range("N2").ClearContents()
$ws.Range("N116").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 850.5625
$ws.Range("I3").Value = 850.5625
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 850.5625
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -736.5625
$ws.Range("H20").Value = 1487.9231
$ws.Range("I20").Value = 1239.375
$ws.Range("J20").Value = 1885.6
$ws.Range("K20").Value = 1239.375
$ws.Range("L20").Value = 1885.6
$ws.Range("M20").Value = -992.375
$ws.Range("N20").Value = -2379.6
$ws.Range("H134").Value = 1465.721
$ws.Range("I134").Value = 795.871
$ws.Range("J134").Value = 3196.1667
$ws.Range("K134").Value = 2387.613
$ws.Range("L134").Value = 9588.500100000001
$ws.Range("M134").Value = 147.3870000000002
$ws.Range("N134").Value = -14658.5001
$ws.Range("N3").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2489.1396
$ws.Range("I31").Value = 892.32355
$ws.Range("J31").Value = 8521.556
$ws.Range("K31").Value = 892.32355
$ws.Range("L31").Value = 8521.556
$ws.Range("M31").Value = -597.32355
$ws.Range("N31").Value = -9111.556
$ws.Range("H34").Value = 2489.1396
$ws.Range("I34").Value = 892.32355
$ws.Range("J34").Value = 8521.556
$ws.Range("K34").Value = 892.32355
$ws.Range("L34").Value = 8521.556
$ws.Range("M34").Value = -690.32355
$ws.Range("N34").Value = -8925.556
$ws.Range("H99").Value = 57685.832
$ws.Range("I99").Value = 144057.42
$ws.Range("J99").Value = 2722.0908
$ws.Range("K99").Value = 144057.42
$ws.Range("L99").Value = 2722.0908
$ws.Range("M99").Value = -142559.42
$ws.Range("N99").Value = -5718.0908
$ws.Range("H105").Value = 814.0769
$ws.Range("I105").Value = 790.3103599999999
$ws.Range("J105").Value = 883
$ws.Range("K105").Value = 790.3103599999999
$ws.Range("L105").Value = 883
$ws.Range("M105").Value = 956.6896400000001
$ws.Range("N105").Value = -4377
$ws.Range("H126").Value = 57685.832
$ws.Range("I126").Value = 144057.42
$ws.Range("J126").Value = 2722.0908
$ws.Range("K126").Value = 432172.26
$ws.Range("L126").Value = 8166.2724
$ws.Range("M126").Value = -429702.26
$ws.Range("N126").Value = -13106.2724
$ws.Range("H134").Value = 2114.6
$ws.Range("I134").Value = 1183.2858
$ws.Range("J134").Value = 3299.9092
$ws.Range("K134").Value = 3549.8574
$ws.Range("L134").Value = 9899.7276
$ws.Range("M134").Value = -1014.8574
$ws.Range("N134").Value = -14969.7276
$ws.Range("H135").Value = 44457.5
$ws.Range("J135").Value = 31554
$ws.Range("L135").Value = 31554
$ws.Range("N135").Value = -41694

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2376.0588
$ws.Range("I69").Value = 633.2857
$ws.Range("J69").Value = 3596
$ws.Range("K69").Value = 1899.8571
$ws.Range("L69").Value = 10788
$ws.Range("M69").Value = -1088.8571
$ws.Range("N69").Value = -12410
$ws.Range("H72").Value = 2376.0588
$ws.Range("I72").Value = 633.2857
$ws.Range("J72").Value = 3596
$ws.Range("K72").Value = 5699.571300000001
$ws.Range("L72").Value = 32364
$ws.Range("M72").Value = -1643.571300000001
$ws.Range("N72").Value = -40476

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2062.1
$ws.Range("I7").Value = 1760.9474
$ws.Range("J7").Value = 2582.2727
$ws.Range("K7").Value = 1760.9474
$ws.Range("L7").Value = 2582.2727
$ws.Range("M7").Value = -1648.9474
$ws.Range("N7").Value = -2806.2727
$ws.Range("H40").Value = 2456.6667
$ws.Range("I40").Value = 2250
$ws.Range("J40").Value = 3025
$ws.Range("K40").Value = 2250
$ws.Range("L40").Value = 3025
$ws.Range("M40").Value = -2114
$ws.Range("N40").Value = -3297
$ws.Range("H61").Value = 1179.4
$ws.Range("I61").Value = 1179.4
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1179.4
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -977.4000000000001
$ws.Range("H82").Value = 1886.8572
$ws.Range("I82").Value = 1634
$ws.Range("J82").Value = 2076.5
$ws.Range("K82").Value = 1634
$ws.Range("L82").Value = 2076.5
$ws.Range("M82").Value = -1273
$ws.Range("N82").Value = -2798.5
$ws.Range("H85").Value = 1886.8572
$ws.Range("I85").Value = 1634
$ws.Range("J85").Value = 2076.5
$ws.Range("K85").Value = 1634
$ws.Range("L85").Value = 2076.5
$ws.Range("M85").Value = -386
$ws.Range("N85").Value = -4572.5
$ws.Range("H113").Value = 1179.4
$ws.Range("I113").Value = 1179.4
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1179.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 990.5999999999999
$ws.Range("H126").Value = 2062.1
$ws.Range("I126").Value = 1760.9474
$ws.Range("J126").Value = 2582.2727
$ws.Range("K126").Value = 5282.8422
$ws.Range("L126").Value = 7746.8181
$ws.Range("M126").Value = -2812.8422
$ws.Range("N126").Value = -12686.8181
$ws.Range("N61").ClearContents()
$ws.Range("N113").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 47503.5
$ws.Range("J46").Value = 47503.5
$ws.Range("L46").Value = 47503.5
$ws.Range("N46").Value = -47965.5
$ws.Range("H134").Value = 47503.5
$ws.Range("J134").Value = 47503.5
$ws.Range("L134").Value = 142510.5
$ws.Range("N134").Value = -147580.5
